# Weekly update: insert a new "Cebollín" price block (Extra/Primera/Segunda/
# Tercera) for the most recent week (2021-10-07) at the top of the dated
# history (row 523), pushing the rest of the table down by 4 rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at 523, shifting existing rows 523:622 down to 527:626
# (this also drags row formatting - e.g. the date style on column D - down
# with it, and grows the sheet's used-range/dimension automatically).
$ws.Range("A523:A526").EntireRow.Insert()

# Row 523 - Extra
$ws.Range("A523").Value = 9
$ws.Range("B523").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C523").Value = "Metropolitana"
$ws.Range("D523").Value = 44476
$ws.Range("E523").Value = 13
$ws.Range("F523").Value = 100112037
$ws.Range("G523").Value = "Cebollín"
$ws.Range("H523").Value = "Sin especificar"
$ws.Range("I523").Value = "Extra"
$ws.Range("J523").Value = 106
$ws.Range("K523").Value = 2500
$ws.Range("L523").Value = 2800
$ws.Range("M523").Value = 2650
$ws.Range("N523").Value = '$/paquete 36 unidades'
$ws.Range("O523").Value = "Región Metropolitana"
$ws.Range("P523").Value = 74
$ws.Range("Q523").Value = 36
$ws.Range("R523").Value = "Hortaliza"

# Row 524 - Primera
$ws.Range("A524").Value = 9
$ws.Range("B524").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C524").Value = "Metropolitana"
$ws.Range("D524").Value = 44476
$ws.Range("E524").Value = 13
$ws.Range("F524").Value = 100112037
$ws.Range("G524").Value = "Cebollín"
$ws.Range("H524").Value = "Sin especificar"
$ws.Range("I524").Value = "Primera"
$ws.Range("J524").Value = 250
$ws.Range("K524").Value = 2100
$ws.Range("L524").Value = 2400
$ws.Range("M524").Value = 2250
$ws.Range("N524").Value = '$/paquete 36 unidades'
$ws.Range("O524").Value = "Región Metropolitana"
$ws.Range("P524").Value = 62
$ws.Range("Q524").Value = 36
$ws.Range("R524").Value = "Hortaliza"

# Row 525 - Segunda
$ws.Range("A525").Value = 9
$ws.Range("B525").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C525").Value = "Metropolitana"
$ws.Range("D525").Value = 44476
$ws.Range("E525").Value = 13
$ws.Range("F525").Value = 100112037
$ws.Range("G525").Value = "Cebollín"
$ws.Range("H525").Value = "Sin especificar"
$ws.Range("I525").Value = "Segunda"
$ws.Range("J525").Value = 160
$ws.Range("K525").Value = 1700
$ws.Range("L525").Value = 2000
$ws.Range("M525").Value = 1850
$ws.Range("N525").Value = '$/paquete 36 unidades'
$ws.Range("O525").Value = "Región Metropolitana"
$ws.Range("P525").Value = 51
$ws.Range("Q525").Value = 36
$ws.Range("R525").Value = "Hortaliza"

# Row 526 - Tercera
$ws.Range("A526").Value = 9
$ws.Range("B526").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C526").Value = "Metropolitana"
$ws.Range("D526").Value = 44476
$ws.Range("E526").Value = 13
$ws.Range("F526").Value = 100112037
$ws.Range("G526").Value = "Cebollín"
$ws.Range("H526").Value = "Sin especificar"
$ws.Range("I526").Value = "Tercera"
$ws.Range("J526").Value = 79
$ws.Range("K526").Value = 1400
$ws.Range("L526").Value = 1600
$ws.Range("M526").Value = 1499
$ws.Range("N526").Value = '$/paquete 36 unidades'
$ws.Range("O526").Value = "Región Metropolitana"
$ws.Range("P526").Value = 42
$ws.Range("Q526").Value = 36
$ws.Range("R526").Value = "Hortaliza"
